$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1299.99
$summary.Range("B4").Value = -0.01
$summary.Range("B5").Value = -0.03
$summary.Range("B6").Value = 7
$summary.Range("B7").Value = 4
$summary.Range("B9").Value = 57.14

# --- Strategy Status sheet (MarketMaking row) -----------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98999999999999
$status.Range("D4").Value = 7
$status.Range("E4").Value = -0.01
$status.Range("F4").Value = -0.01
$status.Range("G4").Value = 57.14

# --- All Trades + MarketMaking sheets: append trade #7 (row 8) -----------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Seed row 8 by copying row 7 first so text-like cells (dates/strings)
    # land as plain text instead of being auto-coerced into dates/numbers
    # when assigned directly via .Value. B2:B7 already hold the literal
    # string "2026-02-17" (same date as the new trade), so leave B8 as-is
    # from the copy instead of re-assigning it (re-assigning a string that
    # looks like a date would get auto-converted to a date serial).
    $ws.Range("A7:Q7").Copy($ws.Range("A8:Q8"))

    $ws.Range("A8").Value = 7
    $ws.Range("C8").Value = "19:44:18"
    $ws.Range("D8").Value = "MarketMaking"
    $ws.Range("E8").Value = "DOWN"
    $ws.Range("F8").Value = 0.93
    $ws.Range("G8").Value = 0.95
    $ws.Range("H8").Value = "CLOSED"
    $ws.Range("I8").Value = 2.1505
    $ws.Range("J8").Value = 0.02
    $ws.Range("K8").Value = 99.98999999999999
    $ws.Range("L8").Value = 0
    $ws.Range("M8").Value = 0
    $ws.Range("N8").Value = 0.6
    $ws.Range("O8").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P8").Value = "early_exit"
    $ws.Range("Q8").Value = 0.15
}
